$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dValues = @{
    2 = "97.304.79"
    3 = "3.718.32"
    5 = "236.85"
    6 = "1.95"
    7 = "657.65"
    11 = "3.716.39"
    13 = "44.85"
    15 = "6.91"
    16 = "4.415.08"
    17 = "96.992.34"
    19 = "3.718.23"
    20 = "13.04"
    21 = "18.75"
    22 = "0.509"
    23 = "524.75"
    24 = "3.47"
    27 = "106.50"
    28 = "0.194"
    29 = "3.920.60"
    30 = "13.52"
    33 = "1.00"
    35 = "1.84"
    36 = "32.62"
    37 = "1.00"
    38 = "643.80"
    39 = "0.593"
    40 = "8.75"
    42 = "0.166"
    43 = "0.498"
    45 = "40.72"
    47 = "0.968"
    50 = "23.62"
    51 = "8.67"
}

$eValues = @{
    2 = "  +0.33%  "
    3 = "  +1.27%  "
    4 = "  +0.06%  "
    5 = "  -1.25%  "
    6 = "  +3.12%  "
    7 = "  +0.29%  "
    8 = "  +2.89%  "
    9 = "  +0.03%  "
    10 = "  -1.74%  "
    11 = "  +1.22%  "
    12 = "  +18.86%  "
    13 = "  -1.57%  "
    14 = "  +0.69%  "
    15 = "  +0.90%  "
    16 = "  +1.29%  "
    17 = "  +0.34%  "
    18 = "  +0.60%  "
    19 = "  +1.63%  "
    20 = "  +2.01%  "
    21 = "  -1.01%  "
    22 = "  -4.28%  "
    23 = "  -1.49%  "
    24 = "  -0.97%  "
    25 = "  +10.89%  "
    26 = "  -4.37%  "
    27 = "  +3.96%  "
    28 = "  +16.08%  "
    29 = "  +1.36%  "
    30 = "  +0.00%  "
    31 = "  +0.25%  "
    32 = "  -1.05%  "
    33 = "  +0.12%  "
    34 = "  +3.43%  "
    35 = "  -2.42%  "
    36 = "  -0.40%  "
    37 = "  -0.12%  "
    38 = "  -1.50%  "
    39 = "  -1.00%  "
    40 = "  -1.80%  "
    42 = "  +1.37%  "
    43 = "  +12.59%  "
    44 = "  -1.82%  "
    45 = "  +4.95%  "
    46 = "  +1.12%  "
    47 = "  +0.43%  "
    48 = "  -0.52%  "
    49 = "  +2.71%  "
    50 = "  -0.05%  "
    51 = "  -1.09%  "
}

foreach ($row in $dValues.Keys) {
    $ws.Range("D$row").Value = $dValues[$row]
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}
